$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Init")

$ws.Range("D5").Value = "A26"
$ws.Range("D6").Value = "B26"
$ws.Range("D7").Value = "C26"
$ws.Range("D8").Value = "G26"
$ws.Range("D9").Value = "H26"
$ws.Range("D10").Value = "I26"
$ws.Range("D11").Value = "J26"
